$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 0.02354566666666667
$ws.Range("H2").Value2 = 0.07063700000000001
$ws.Range("I2").Value2 = 0.002815555392485919
$ws.Range("J2").Value2 = 0.002815555392485918
$ws.Range("K2").Value2 = 3
$ws.Range("L2").Value2 = 1
$ws.Range("M2").Value2 = 0.9218943333333334
$ws.Range("N2").Value2 = 2.765683
$ws.Range("O2").Value2 = 0.04284983107934069
$ws.Range("P2").Value2 = 0.04284983107934069
$ws.Range("Q2").Value2 = 0.02170661667455556
$ws.Range("R2").Value2 = 0.195359550071
$ws.Range("S2").Value2 = 0.0001206460729625484
$ws.Range("T2").Value2 = 0.0001206460729625484

$ws.Range("G3").Value2 = 0.02354566666666667
$ws.Range("H3").Value2 = 0.07063700000000001
$ws.Range("I3").Value2 = 0.002815555392485919
$ws.Range("J3").Value2 = 0.002815555392485918
$ws.Range("O3").Value2 = 0.3860239622463043
$ws.Range("P3").Value2 = 0.3860239622463043
$ws.Range("Q3").Value2 = 0.1955497598148889
$ws.Range("R3").Value2 = 1.759947838334
$ws.Range("S3").Value2 = 0.001086871848531363
$ws.Range("T3").Value2 = 0.001086871848531363

$ws.Range("G4").Value2 = 0.02354566666666667
$ws.Range("H4").Value2 = 0.07063700000000001
$ws.Range("I4").Value2 = 0.002815555392485919
$ws.Range("J4").Value2 = 0.002815555392485918
$ws.Range("O4").Value2 = 0.571126206674355
$ws.Range("P4").Value2 = 0.571126206674355
$ws.Range("Q4").Value2 = 0.2893177715944445
$ws.Range("R4").Value2 = 2.60385994435
$ws.Range("S4").Value2 = 0.001608037470992007
$ws.Range("T4").Value2 = 0.001608037470992007

$ws.Range("I5").Value2 = 0.9868456480383168
$ws.Range("J5").Value2 = 0.9868456480383166
$ws.Range("K5").Value2 = 3
$ws.Range("L5").Value2 = 1
$ws.Range("M5").Value2 = 0.9218943333333334
$ws.Range("N5").Value2 = 2.765683
$ws.Range("O5").Value2 = 0.04284983107934069
$ws.Range("P5").Value2 = 0.04284983107934069
$ws.Range("Q5").Value2 = 7.608118901190556
$ws.Range("R5").Value2 = 68.473070110715
$ws.Range("S5").Value2 = 0.04228616931982437
$ws.Range("T5").Value2 = 0.04228616931982437

$ws.Range("I6").Value2 = 0.9868456480383168
$ws.Range("J6").Value2 = 0.9868456480383166
$ws.Range("O6").Value2 = 0.3860239622463043
$ws.Range("P6").Value2 = 0.3860239622463043
$ws.Range("S6").Value2 = 0.3809460671812728
$ws.Range("T6").Value2 = 0.3809460671812728

$ws.Range("I7").Value2 = 0.9868456480383168
$ws.Range("J7").Value2 = 0.9868456480383166
$ws.Range("O7").Value2 = 0.571126206674355
$ws.Range("P7").Value2 = 0.571126206674355
$ws.Range("S7").Value2 = 0.5636134115372194
$ws.Range("T7").Value2 = 0.5636134115372194

$ws.Range("I8").Value2 = 0.0103387965691973
$ws.Range("J8").Value2 = 0.0103387965691973
$ws.Range("K8").Value2 = 3
$ws.Range("L8").Value2 = 1
$ws.Range("M8").Value2 = 0.9218943333333334
$ws.Range("N8").Value2 = 2.765683
$ws.Range("O8").Value2 = 0.04284983107934069
$ws.Range("P8").Value2 = 0.04284983107934069
$ws.Range("Q8").Value2 = 0.0797072913581111
$ws.Range("R8").Value2 = 0.7173656222229999
$ws.Range("S8").Value2 = 0.0004430156865537715
$ws.Range("T8").Value2 = 0.0004430156865537715

$ws.Range("I9").Value2 = 0.0103387965691973
$ws.Range("J9").Value2 = 0.0103387965691973
$ws.Range("O9").Value2 = 0.3860239622463043
$ws.Range("P9").Value2 = 0.3860239622463043
$ws.Range("S9").Value2 = 0.00399102321650004
$ws.Range("T9").Value2 = 0.00399102321650004

$ws.Range("I10").Value2 = 0.0103387965691973
$ws.Range("J10").Value2 = 0.0103387965691973
$ws.Range("O10").Value2 = 0.571126206674355
$ws.Range("P10").Value2 = 0.571126206674355
$ws.Range("S10").Value2 = 0.005904757666143491
$ws.Range("T10").Value2 = 0.005904757666143491
